{"js": "// Replace the date line and each two-digit \u00f7 one-digit division problem\n// with the new values from the commit, using exact, unique text matches.\nconst replacements = [\n  [\"2025-10-22 Wednesday\", \"2025-10-23 Thursday\"],\n  [\"76\u00f77=10, 6\", \"67\u00f72=33, 1\"],\n  [\"50\u00f76=8, 2\", \"93\u00f75=18, 3\"],\n  [\"82\u00f77=11, 5\", \"87\u00f77=12, 3\"],\n  [\"82\u00f76=13, 4\", \"35\u00f75=7, 0\"],\n  [\"46\u00f72=23, 0\", \"71\u00f76=11, 5\"],\n  [\"31\u00f79=3, 4\", \"48\u00f72=24, 0\"],\n  [\"97\u00f77=13, 6\", \"98\u00f76=16, 2\"],\n  [\"83\u00f75=16, 3\", \"53\u00f75=10, 3\"],\n  [\"88\u00f78=11, 0\", \"18\u00f78=2, 2\"],\n  [\"65\u00f73=21, 2\", \"17\u00f77=2, 3\"],\n  [\"60\u00f78=7, 4\", \"37\u00f75=7, 2\"],\n  [\"26\u00f72=13, 0\", \"86\u00f77=12, 2\"],\n  [\"16\u00f76=2, 4\", \"21\u00f79=2, 3\"],\n  [\"15\u00f76=2, 3\", \"13\u00f78=1, 5\"],\n  [\"33\u00f78=4, 1\", \"50\u00f76=8, 2\"],\n  [\"75\u00f73=25, 0\", \"71\u00f77=10, 1\"],\n  [\"69\u00f75=13, 4\", \"33\u00f79=3, 6\"],\n  [\"37\u00f74=9, 1\", \"80\u00f75=16, 0\"],\n  [\"74\u00f76=12, 2\", \"96\u00f74=24, 0\"],\n  [\"61\u00f73=20, 1\", \"26\u00f79=2, 8\"],\n  [\"11\u00f77=1, 4\", \"98\u00f72=49, 0\"],\n  [\"50\u00f77=7, 1\", \"89\u00f75=17, 4\"],\n  [\"86\u00f72=43, 0\", \"81\u00f79=9, 0\"],\n  [\"43\u00f75=8, 3\", \"26\u00f77=3, 5\"],\n  [\"79\u00f77=11, 2\", \"44\u00f75=8, 4\"],\n];\n\nconst body = context.document.body;\n\nfor (const [oldText, newText] of replacements) {\n  const results = body.search(oldText, { matchCase: true, matchWholeWord: false });\n  results.load(\"items\");\n  await context.sync();\n\n  if (results.items.length === 0) {\n    throw new Error(`Text not found: ${oldText}`);\n  }\n\n  for (const range of results.items) {\n    range.insertText(newText, \"Replace\");\n  }\n}\n\nawait context.sync();\n", "ps1": "# Update the worksheet date and every two-digit \u00f7 one-digit division\n# problem to the new values, matching exact text via Find & Replace.\n$d = $word.ActiveDocument\n\n$pairs = @(\n    @(\"2025-10-22 Wednesday\", \"2025-10-23 Thursday\"),\n    @(\"76\u00f77=10, 6\", \"67\u00f72=33, 1\"),\n    @(\"50\u00f76=8, 2\", \"93\u00f75=18, 3\"),\n    @(\"82\u00f77=11, 5\", \"87\u00f77=12, 3\"),\n    @(\"82\u00f76=13, 4\", \"35\u00f75=7, 0\"),\n    @(\"46\u00f72=23, 0\", \"71\u00f76=11, 5\"),\n    @(\"31\u00f79=3, 4\", \"48\u00f72=24, 0\"),\n    @(\"97\u00f77=13, 6\", \"98\u00f76=16, 2\"),\n    @(\"83\u00f75=16, 3\", \"53\u00f75=10, 3\"),\n    @(\"88\u00f78=11, 0\", \"18\u00f78=2, 2\"),\n    @(\"65\u00f73=21, 2\", \"17\u00f77=2, 3\"),\n    @(\"60\u00f78=7, 4\", \"37\u00f75=7, 2\"),\n    @(\"26\u00f72=13, 0\", \"86\u00f77=12, 2\"),\n    @(\"16\u00f76=2, 4\", \"21\u00f79=2, 3\"),\n    @(\"15\u00f76=2, 3\", \"13\u00f78=1, 5\"),\n    @(\"33\u00f78=4, 1\", \"50\u00f76=8, 2\"),\n    @(\"75\u00f73=25, 0\", \"71\u00f77=10, 1\"),\n    @(\"69\u00f75=13, 4\", \"33\u00f79=3, 6\"),\n    @(\"37\u00f74=9, 1\", \"80\u00f75=16, 0\"),\n    @(\"74\u00f76=12, 2\", \"96\u00f74=24, 0\"),\n    @(\"61\u00f73=20, 1\", \"26\u00f79=2, 8\"),\n    @(\"11\u00f77=1, 4\", \"98\u00f72=49, 0\"),\n    @(\"50\u00f77=7, 1\", \"89\u00f75=17, 4\"),\n    @(\"86\u00f72=43, 0\", \"81\u00f79=9, 0\"),\n    @(\"43\u00f75=8, 3\", \"26\u00f77=3, 5\"),\n    @(\"79\u00f77=11, 2\", \"44\u00f75=8, 4\")\n)\n\nforeach ($pair in $pairs) {\n    $oldText = $pair[0]\n    $newText = $pair[1]\n\n    $find = $d.Content.Find\n    $find.ClearFormatting()\n    $find.Replacement.ClearFormatting()\n    $found = $find.Execute($oldText, $false, $false, $false, $false, $false, $true, 1, $false, $newText, 2)\n\n    if (-not $found) {\n        throw \"Text not found: $oldText\"\n    }\n}\n"}
